$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"

$ws.Range("A2").Value = 12
$ws.Range("B2").Value = 44
$ws.Range("C2").Value = 90

$ws.Range("C2").Select()
